$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Lower Right Cell" references for the scenario index blocks
# from row 18 to row 19 (new scenario 19 added).
$ws.Range("D5").Value = "A19"
$ws.Range("D6").Value = "B19"
$ws.Range("D7").Value = "C19"
$ws.Range("D8").Value = "G19"
$ws.Range("D9").Value = "H19"
$ws.Range("D10").Value = "I19"
$ws.Range("D11").Value = "J19"

$ws.Range("D11").Select()
